$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# Row 7 was Arshia Mehta (with email); becomes Danielle Aldrett (no email)
$ws.Range("A7").Value = "Danielle"
$ws.Range("C7").Value = "Aldrett"
$ws.Range("E7").ClearContents()

# Row 8 was Danielle Aldrett (no email); becomes Arshia Mehta (with email)
$ws.Range("A8").Value = "Arshia"
$ws.Range("C8").Value = "Mehta"
$ws.Range("E8").Value = "amehta3@wellesley.edu"

# Row 9 had surname/givenname swapped; fix it: A9 givenName=Rachel, C9 surName=Shrives
$ws.Range("A9").Value = "Rachel"
$ws.Range("C9").Value = "Shrives"

# Update selection on this sheet
$ws.Activate()
$ws.Range("E14").Select()
